$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "season_ending_year_x"
$ws.Range("O1").Value = "season_ending_year_y"

$ws.Range("AY1").Value = "calendar_year"
$ws.Range("AY1").Font.Bold = $true
$ws.Range("AY1").HorizontalAlignment = -4108
$ws.Range("AY1").VerticalAlignment = -4160
$ws.Range("AY1").Borders.LineStyle = 1

$ws.Range("AY2").Value = 1975
$ws.Range("AY3").Value = 1974
$ws.Range("AY4").Value = 1974
$ws.Range("AY5").Value = 1974
$ws.Range("AY6").Value = 1973
$ws.Range("AY7").Value = 1972

$ws.Range("Q2").Value = 1953
$ws.Range("Q6").Value = 1952
$ws.Range("Q7").Value = 1950
